# Refine the data assembling code, generate dataset ready for modelling in NONMEM
#
# The "WT" sheet lists three biological-replicate blocks (Erlotinib /
# ATUX 4,5 uM / DMSO), each with a header/formula row followed by three
# "Tech rep" data rows. This edit inserts one new blank row above the
# third block (before row 12, the "Erlotinib" label that starts the
# third replicate section) so the block headers/labels get their own
# row separate from the formula row beneath them - mirroring the layout
# already used for the first two blocks - and everything below shifts
# down by one row accordingly.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("WT")
$ws.Activate()

# Insert a new row at row 12; existing rows 12:21 (and their formulas)
# shift down to 13:22, so A12's lone "Erlotinib" label becomes A13, the
# ATUX 4,5 uM header/formula row becomes row 14, its three tech-rep rows
# become 15:17, the DMSO header/formula row becomes row 19, and its
# three tech-rep rows become 20:22.
$ws.Rows.Item(12).Insert()

# Matches the author's final cursor position recorded in the workbook.
$ws.Range("E12").Select()
